# Apply ticker-symbol updates to Sheet1 and shrink the used range from
# A1:F20 down to A1:F15 by deleting the now-unused trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update / clear cell contents for rows 2-15 (columns B:F) ---------

# Row 2
$ws.Range("B2").Value = "NSE:ANGELONE"
$ws.Range("C2").Value = "NSE:AARTECH"
$ws.Range("D2").Value = "NSE:GMRINFRA"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "NSE:BPCL"

# Row 3
$ws.Range("B3").Value = "NSE:AXISILVER"
$ws.Range("C3").Value = "NSE:AJANTPHARM"
$ws.Range("D3").Value = "NSE:IOC"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "NSE:HINDPETRO"

# Row 4
$ws.Range("B4").Value = "NSE:BSOFT"
$ws.Range("C4").Value = "NSE:FAZE3Q"
$ws.Range("D4").Value = "NSE:RAMCOCEM"
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""

# Row 5
$ws.Range("B5").Value = "NSE:CAMS"
$ws.Range("C5").Value = "NSE:RATNAMANI"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""

# Row 6
$ws.Range("B6").Value = "NSE:CHEVIOT"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""

# Row 7
$ws.Range("B7").Value = "NSE:ESG"
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""

# Row 8
$ws.Range("B8").Value = "NSE:GRAVITA"
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""

# Row 9
$ws.Range("B9").Value = "NSE:HINDPETRO"
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = ""

# Row 10
$ws.Range("B10").Value = "NSE:KELLTONTEC"
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = ""

# Row 11
$ws.Range("B11").Value = "NSE:LOVABLE"
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""

# Row 12
$ws.Range("B12").Value = "NSE:NEXT50"
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""

# Row 13
$ws.Range("B13").Value = "NSE:PDMJEPAPER"
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = ""

# Row 14
$ws.Range("B14").Value = "NSE:QUICKHEAL"
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = ""

# Row 15
$ws.Range("B15").Value = "NSE:RAMRAT"
$ws.Range("C15").Value = ""
$ws.Range("D15").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("F15").Value = ""

# --- Remove the now-unused rows 16-20 (shrinks dimension to A1:F15) ---

$ws.Range("A16:F20").Delete()
